# Auto-generated Excel COM-interop script to update the cryptos list
# (price/volume refresh + two coin-row swaps), per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.933.61"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "1.645.61"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.62"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.41"
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0872"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("D12").Value = "1.879.77"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "1.644.23"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.60"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "27.947.14"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.26"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.55"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.75"
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "1.442.66"
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +3.32%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0169"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.932"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.559"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.14"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +5.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.42"
$ws.Range("E46").Value = "  +3.30%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "1.787.64"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "89.08"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("E51").Value = "  +0.20%  "
